# "Add files via upload" - Missing road data.xlsx
# Renames the sheet, refreshes the header labels to note the data source
# (disaster.ninja), and updates the sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" -> "Missing_data"
$ws.Name = "Missing_data"

# Refresh header row (row 1) text - same column meanings as before, but
# most labels now note they were extracted from disaster.ninja.
$ws.Range("A1").Value = "Country"
$ws.Range("B1").Value = "Population - extracted from disaster.ninja"
$ws.Range("C1").Value = "Populated area (km2) - extracted from disaster.ninja"
$ws.Range("D1").Value = "Populated area with no road count (km2) - extracted from disaster.ninja"
$ws.Range("E1").Value = "Populated area with no road information (%)"
$ws.Range("F1").Value = "Total land area (km2) - extracted from disaster.ninja"

# Update the view: scroll so row 19 is at the top, and select D1.
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D1").Select()
